# Apply the "testData.xlsx" revision:
#  - PostDetails sheet: drop the merged C1:D1 "STPS" header block (columns C/D)
#    and rename the Plant1-5 / "Details entered successfully for plantN" data
#    to Master1-5 / "Details entered successfully for MasterN".
#  - CreateSTP_Mandatory sheet: rename the Nimda*/TEst*/Test4 data to the new
#    Kind* naming scheme and append 3 more data rows (rows 7-9), re-using the
#    "valid Mandatory Fields" success message and "sharath sethu" values.
#  - Make CreateSTP_Mandatory the active sheet/tab (it was PostDetails before).

$wb = $excel.ActiveWorkbook

$wsPost = $wb.Worksheets.Item("PostDetails")
$wsMand = $wb.Worksheets.Item("CreateSTP_Mandatory")

# ---------------------------------------------------------------------------
# PostDetails: remove the old merged header block in C1:D1 and the stray
# data in C2:D2, then rename the Plant* rows to Master*.
# ---------------------------------------------------------------------------
$wsPost.Range("C1:D1").UnMerge()
$wsPost.Columns.Item(3).Delete()
$wsPost.Columns.Item(3).Delete()

$wsPost.Range("A2").Value = "Master1"
$wsPost.Range("A3").Value = "Master2"
$wsPost.Range("A4").Value = "Master3"
$wsPost.Range("A5").Value = "Master4"
$wsPost.Range("A6").Value = "Master5"
$wsPost.Range("B2").Value = "Details entered successfully for Master1"
$wsPost.Range("B3").Value = "Details entered successfully for Master2"
$wsPost.Range("B4").Value = "Details entered successfully for Master3"
$wsPost.Range("B5").Value = "Details entered successfully for Master4"
$wsPost.Range("B6").Value = "Details entered successfully for Master5"

# ---------------------------------------------------------------------------
# CreateSTP_Mandatory: rename the existing 5 data rows from Nimda*/TEst*
# naming to Kind*, and append 3 new rows (7-9) with the next Kind* values.
# ---------------------------------------------------------------------------
$successMsg = 'This is a valid Mandatory Fields scenario. It is supposed give a toast message "STP created successfully!"'
$leadContact = "sharath sethu"

# Existing rows 2-6: rewrite column A, then column B (column-major order,
# matching how the shared-string table grows in the source revision).
$wsMand.Range("A2").Value = "Kind1"
$wsMand.Range("A3").Value = "Kind2"
$wsMand.Range("A4").Value = "Kind3"
$wsMand.Range("A5").Value = "Kind4"
$wsMand.Range("A6").Value = "Kind5"

$wsMand.Range("B2").Value = "Kind123"
$wsMand.Range("B3").Value = "Kind124"
$wsMand.Range("B4").Value = "Kind125"
$wsMand.Range("B5").Value = "Kind126"
$wsMand.Range("B6").Value = "Kind127"

$wsMand.Range("C2").Value = $successMsg
$wsMand.Range("C3").Value = $successMsg
$wsMand.Range("C4").Value = $successMsg
$wsMand.Range("C5").Value = $successMsg
$wsMand.Range("C6").Value = $successMsg

$wsMand.Range("D2").Value = $leadContact
$wsMand.Range("D3").Value = $leadContact
$wsMand.Range("D4").Value = $leadContact
$wsMand.Range("D5").Value = $leadContact
$wsMand.Range("D6").Value = $leadContact

# New rows 7-9: append column A, then column B, C, D.
$wsMand.Range("A7").Value = "Kind6"
$wsMand.Range("A8").Value = "Kind7"
$wsMand.Range("A9").Value = "Kind8"

$wsMand.Range("B7").Value = "Kind128"
$wsMand.Range("B8").Value = "Kind129"
$wsMand.Range("B9").Value = "Kind130"

$wsMand.Range("C7").Value = $successMsg
$wsMand.Range("C8").Value = $successMsg
$wsMand.Range("C9").Value = $successMsg

$wsMand.Range("D7").Value = $leadContact
$wsMand.Range("D8").Value = $leadContact
$wsMand.Range("D9").Value = $leadContact

# ---------------------------------------------------------------------------
# Selections / active tab: CreateSTP_Mandatory becomes the selected/active
# sheet, PostDetails loses its selection highlight.
# ---------------------------------------------------------------------------
$wsPost.Range("B2:B6").Select()

$wsMand.Activate()
$wsMand.Range("D4:D9").Select()
